# Backlog.xlsx: add a new pending task row to the "Hoja1" backlog sheet.
#   A55 = "agregar filto a frontend de remitos"
#   B55 = "no comenzado"   (same status used by the other newly added rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A55").Value = "agregar filto a frontend de remitos"
$ws.Range("B55").Value = "no comenzado"

# Mirror the author's on-screen state: the new rows selected, scrolled so
# row 29 is at the top of the viewport.
$ws.Activate()
$ws.Range("A55:A56").Select()
$excel.ActiveWindow.ScrollRow = 29
